# This script applies the scheduled-runner profit-recalculation update to all
# 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) in the workbook.
# Each target cell's numeric value is set to match refreshed market-board data;
# a few cells that no longer have a profit/loss figure are cleared, and a couple
# of previously-empty cells gain a new figure.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 118.375  # H8: 105.333336 -> 118.375
$ws.Cells.Item(8, 9).Value = 118.375  # I8: 105.333336 -> 118.375
$ws.Cells.Item(8, 11).Value = 355.125  # K8: 316.000008 -> 355.125
$ws.Cells.Item(8, 13).Value = -216.125  # M8: -177.000008 -> -216.125
$ws.Cells.Item(17, 8).Value = 285514.34  # H17: 302304.12 -> 285514.34
$ws.Cells.Item(17, 10).Value = 285514.34  # J17: 302304.12 -> 285514.34
$ws.Cells.Item(17, 12).Value = 856543.02  # L17: 906912.36 -> 856543.02
$ws.Cells.Item(17, 14).Value = -856879.02  # N17: -907248.36 -> -856879.02
$ws.Cells.Item(32, 8).Value = 1651  # H32: 1450.8572 -> 1651
$ws.Cells.Item(32, 9).Value = 1250  # I32: 1117.6 -> 1250
$ws.Cells.Item(32, 10).Value = 1811.4  # J32: 1636 -> 1811.4
$ws.Cells.Item(32, 11).Value = 1250  # K32: 1117.6 -> 1250
$ws.Cells.Item(32, 12).Value = 1811.4  # L32: 1636 -> 1811.4
$ws.Cells.Item(32, 13).Value = -924  # M32: -791.5999999999999 -> -924
$ws.Cells.Item(32, 14).Value = -2463.4  # N32: -2288 -> -2463.4
$ws.Cells.Item(76, 8).Value = 55002800  # H76: 52383760 -> 55002800
$ws.Cells.Item(76, 9).Value = 55002800  # I76: 52383760 -> 55002800
$ws.Cells.Item(76, 11).Value = 55002800  # K76: 52383760 -> 55002800
$ws.Cells.Item(76, 13).Value = -55002485  # M76: -52383445 -> -55002485
$ws.Cells.Item(79, 8).Value = 55002800  # H79: 52383760 -> 55002800
$ws.Cells.Item(79, 9).Value = 55002800  # I79: 52383760 -> 55002800
$ws.Cells.Item(79, 11).Value = 55002800  # K79: 52383760 -> 55002800
$ws.Cells.Item(79, 13).Value = -55001708  # M79: -52382668 -> -55001708
$ws.Cells.Item(132, 8).Value = 1472.1351  # H132: 1579.0714 -> 1472.1351
$ws.Cells.Item(132, 9).Value = 1278.409  # I132: 1384.3064 -> 1278.409
$ws.Cells.Item(132, 10).Value = 3070.375  # J132: 3088.5 -> 3070.375
$ws.Cells.Item(132, 11).Value = 3835.227  # K132: 4152.9192 -> 3835.227
$ws.Cells.Item(132, 12).Value = 9211.125  # L132: 9265.5 -> 9211.125
$ws.Cells.Item(132, 13).Value = -1305.227  # M132: -1622.9192 -> -1305.227
$ws.Cells.Item(132, 14).Value = -14271.125  # N132: -14325.5 -> -14271.125
$ws.Cells.Item(137, 8).Value = 1093.1333  # H137: 1057.3438 -> 1093.1333
$ws.Cells.Item(137, 9).Value = 772.86365  # I137: 751.8333 -> 772.86365
$ws.Cells.Item(137, 11).Value = 2318.59095  # K137: 2255.4999 -> 2318.59095
$ws.Cells.Item(137, 13).Value = 231.4090500000002  # M137: 294.5001000000002 -> 231.4090500000002
$ws.Cells.Item(138, 8).Value = 2721.5305  # H138: 2745.23 -> 2721.5305
$ws.Cells.Item(138, 9).Value = 662.36365  # I138: 656.69696 -> 662.36365
$ws.Cells.Item(138, 10).Value = 3766.9539  # J138: 3773.9104 -> 3766.9539
$ws.Cells.Item(138, 11).Value = 1987.09095  # K138: 1970.09088 -> 1987.09095
$ws.Cells.Item(138, 12).Value = 11300.8617  # L138: 11321.7312 -> 11300.8617
$ws.Cells.Item(138, 13).Value = 3152.90905  # M138: 3169.90912 -> 3152.90905
$ws.Cells.Item(138, 14).Value = -21580.8617  # N138: -21601.7312 -> -21580.8617
$ws.Cells.Item(139, 8).Value = 54552  # H139: 55000 -> 54552
$ws.Cells.Item(139, 10).Value = 54552  # J139: 55000 -> 54552
$ws.Cells.Item(139, 12).Value = 54552  # L139: 55000 -> 54552
$ws.Cells.Item(139, 14).Value = -64832  # N139: -65280 -> -64832

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6867.4736  # H32: 8553.333000000001 -> 6867.4736
$ws.Cells.Item(32, 9).Value = 3121.1045  # I32: 4490 -> 3121.1045
$ws.Cells.Item(32, 10).Value = 34757.11  # J32: 44542.855 -> 34757.11
$ws.Cells.Item(32, 11).Value = 3121.1045  # K32: 4490 -> 3121.1045
$ws.Cells.Item(32, 12).Value = 34757.11  # L32: 44542.855 -> 34757.11
$ws.Cells.Item(32, 13).Value = -2834.1045  # M32: -4203 -> -2834.1045
$ws.Cells.Item(32, 14).Value = -35331.11  # N32: -45116.855 -> -35331.11
$ws.Cells.Item(61, 8).Value = 847.1111  # H61: 746.4400000000001 -> 847.1111
$ws.Cells.Item(61, 9).Value = 788.93335  # I61: 693.34784 -> 788.93335
$ws.Cells.Item(61, 10).Value = 1138  # J61: 1357 -> 1138
$ws.Cells.Item(61, 11).Value = 788.93335  # K61: 693.34784 -> 788.93335
$ws.Cells.Item(61, 12).Value = 1138  # L61: 1357 -> 1138
$ws.Cells.Item(61, 13).Value = -576.93335  # M61: -481.34784 -> -576.93335
$ws.Cells.Item(61, 14).Value = -1562  # N61: -1781 -> -1562
$ws.Cells.Item(74, 8).Value = 1021.5909  # H74: 812.3611 -> 1021.5909
$ws.Cells.Item(74, 9).Value = 972.3684  # I74: 758.7241 -> 972.3684
$ws.Cells.Item(74, 10).Value = 1333.3334  # J74: 1034.5714 -> 1333.3334
$ws.Cells.Item(74, 11).Value = 972.3684  # K74: 758.7241 -> 972.3684
$ws.Cells.Item(74, 12).Value = 1333.3334  # L74: 1034.5714 -> 1333.3334
$ws.Cells.Item(74, 13).Value = -98.36839999999995  # M74: 115.2759 -> -98.36839999999995
$ws.Cells.Item(74, 14).Value = -3081.3334  # N74: -2782.5714 -> -3081.3334
$ws.Cells.Item(77, 8).Value = 1021.5909  # H77: 812.3611 -> 1021.5909
$ws.Cells.Item(77, 9).Value = 972.3684  # I77: 758.7241 -> 972.3684
$ws.Cells.Item(77, 10).Value = 1333.3334  # J77: 1034.5714 -> 1333.3334
$ws.Cells.Item(77, 11).Value = 4861.842  # K77: 3793.6205 -> 4861.842
$ws.Cells.Item(77, 12).Value = 6666.666999999999  # L77: 5172.857 -> 6666.666999999999
$ws.Cells.Item(77, 13).Value = -493.8419999999996  # M77: 574.3795 -> -493.8419999999996
$ws.Cells.Item(77, 14).Value = -15402.667  # N77: -13908.857 -> -15402.667
$ws.Cells.Item(101, 8).Value = 22000  # H101: 20750 -> 22000
$ws.Cells.Item(101, 10).Value = 22000  # J101: 20750 -> 22000
$ws.Cells.Item(101, 12).Value = 22000  # L101: 20750 -> 22000
$ws.Cells.Item(101, 14).Value = -28490  # N101: -27240 -> -28490
$ws.Cells.Item(132, 8).Value = 1312.5555  # H132: 1558.2858 -> 1312.5555
$ws.Cells.Item(132, 9).Value = 840.78125  # I132: 1017.2727 -> 840.78125
$ws.Cells.Item(132, 11).Value = 2522.34375  # K132: 3051.8181 -> 2522.34375
$ws.Cells.Item(132, 13).Value = 7.65625  # M132: -521.8181 -> 7.65625
$ws.Cells.Item(136, 8).Value = 847.1111  # H136: 746.4400000000001 -> 847.1111
$ws.Cells.Item(136, 9).Value = 788.93335  # I136: 693.34784 -> 788.93335
$ws.Cells.Item(136, 10).Value = 1138  # J136: 1357 -> 1138
$ws.Cells.Item(136, 11).Value = 2366.80005  # K136: 2080.04352 -> 2366.80005
$ws.Cells.Item(136, 12).Value = 3414  # L136: 4071 -> 3414
$ws.Cells.Item(136, 13).Value = 183.1999500000002  # M136: 469.9564799999998 -> 183.1999500000002
$ws.Cells.Item(136, 14).Value = -8514  # N136: -9171 -> -8514
$ws.Cells.Item(137, 8).Value = 73400  # H137: 65275 -> 73400
$ws.Cells.Item(137, 10).Value = 73400  # J137: 65275 -> 73400
$ws.Cells.Item(137, 12).Value = 73400  # L137: 65275 -> 73400
$ws.Cells.Item(137, 14).Value = -83600  # N137: -75475 -> -83600
$ws.Cells.Item(138, 8).Value = 45429  # H138: 0 -> 45429
$ws.Cells.Item(138, 10).Value = 45429  # J138: 0 -> 45429
$ws.Cells.Item(138, 12).Value = 45429  # L138: 0 -> 45429
$ws.Cells.Item(138, 14).Value = -55709  # N138: None -> -55709

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 7459.3335  # H105: 7602 -> 7459.3335
$ws.Cells.Item(105, 9).Value = 8498.888999999999  # I105: 9287.5 -> 8498.888999999999
$ws.Cells.Item(105, 10).Value = 5900  # J105: 5675.7144 -> 5900
$ws.Cells.Item(105, 11).Value = 8498.888999999999  # K105: 9287.5 -> 8498.888999999999
$ws.Cells.Item(105, 12).Value = 5900  # L105: 5675.7144 -> 5900
$ws.Cells.Item(105, 13).Value = -6751.888999999999  # M105: -7540.5 -> -6751.888999999999
$ws.Cells.Item(105, 14).Value = -9394  # N105: -9169.714400000001 -> -9394
$ws.Cells.Item(134, 8).Value = 70700.34  # H134: 85634.086 -> 70700.34
$ws.Cells.Item(134, 9).Value = 1704.2307  # I134: 2624.375 -> 1704.2307
$ws.Cells.Item(134, 10).Value = 668666.7  # J134: 251653.5 -> 668666.7
$ws.Cells.Item(134, 11).Value = 5112.6921  # K134: 7873.125 -> 5112.6921
$ws.Cells.Item(134, 12).Value = 2006000.1  # L134: 754960.5 -> 2006000.1
$ws.Cells.Item(134, 13).Value = -2577.6921  # M134: -5338.125 -> -2577.6921
$ws.Cells.Item(134, 14).Value = -2011070.1  # N134: -760030.5 -> -2011070.1
$ws.Cells.Item(135, 8).Value = 60639.5  # H135: 70779 -> 60639.5
$ws.Cells.Item(135, 10).Value = 60639.5  # J135: 70779 -> 60639.5
$ws.Cells.Item(135, 12).Value = 60639.5  # L135: 70779 -> 60639.5
$ws.Cells.Item(135, 14).Value = -70779.5  # N135: -80919 -> -70779.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1406.8536  # H31: 1968.1277 -> 1406.8536
$ws.Cells.Item(31, 9).Value = 1367.025  # I31: 1945.6957 -> 1367.025
$ws.Cells.Item(31, 11).Value = 1367.025  # K31: 1945.6957 -> 1367.025
$ws.Cells.Item(31, 13).Value = -1072.025  # M31: -1650.6957 -> -1072.025
$ws.Cells.Item(34, 8).Value = 1406.8536  # H34: 1968.1277 -> 1406.8536
$ws.Cells.Item(34, 9).Value = 1367.025  # I34: 1945.6957 -> 1367.025
$ws.Cells.Item(34, 11).Value = 1367.025  # K34: 1945.6957 -> 1367.025
$ws.Cells.Item(34, 13).Value = -1165.025  # M34: -1743.6957 -> -1165.025
$ws.Cells.Item(58, 8).Value = 3789.7837  # H58: 5229.577 -> 3789.7837
$ws.Cells.Item(58, 9).Value = 1193.5454  # I58: 1681.4286 -> 1193.5454
$ws.Cells.Item(58, 10).Value = 7597.6  # J58: 9369.083000000001 -> 7597.6
$ws.Cells.Item(58, 11).Value = 1193.5454  # K58: 1681.4286 -> 1193.5454
$ws.Cells.Item(58, 12).Value = 7597.6  # L58: 9369.083000000001 -> 7597.6
$ws.Cells.Item(58, 13).Value = -990.5454  # M58: -1478.4286 -> -990.5454
$ws.Cells.Item(58, 14).Value = -8003.6  # N58: -9775.083000000001 -> -8003.6
$ws.Cells.Item(99, 8).Value = 2176.7  # H99: 2362.6775 -> 2176.7
$ws.Cells.Item(99, 9).Value = 1855.9231  # I99: 1966.45 -> 1855.9231
$ws.Cells.Item(99, 10).Value = 2772.4285  # J99: 3083.0908 -> 2772.4285
$ws.Cells.Item(99, 11).Value = 1855.9231  # K99: 1966.45 -> 1855.9231
$ws.Cells.Item(99, 12).Value = 2772.4285  # L99: 3083.0908 -> 2772.4285
$ws.Cells.Item(99, 13).Value = -357.9231  # M99: -468.45 -> -357.9231
$ws.Cells.Item(99, 14).Value = -5768.4285  # N99: -6079.0908 -> -5768.4285
$ws.Cells.Item(126, 8).Value = 2176.7  # H126: 2362.6775 -> 2176.7
$ws.Cells.Item(126, 9).Value = 1855.9231  # I126: 1966.45 -> 1855.9231
$ws.Cells.Item(126, 10).Value = 2772.4285  # J126: 3083.0908 -> 2772.4285
$ws.Cells.Item(126, 11).Value = 5567.7693  # K126: 5899.35 -> 5567.7693
$ws.Cells.Item(126, 12).Value = 8317.2855  # L126: 9249.2724 -> 8317.2855
$ws.Cells.Item(126, 13).Value = -3097.7693  # M126: -3429.35 -> -3097.7693
$ws.Cells.Item(126, 14).Value = -13257.2855  # N126: -14189.2724 -> -13257.2855
$ws.Cells.Item(132, 8).Value = 1280.6888  # H132: 2304.625 -> 1280.6888
$ws.Cells.Item(132, 9).Value = 747.7742  # I132: 1158.4445 -> 747.7742
$ws.Cells.Item(132, 10).Value = 2460.7144  # J132: 3778.2856 -> 2460.7144
$ws.Cells.Item(132, 11).Value = 2243.3226  # K132: 3475.3335 -> 2243.3226
$ws.Cells.Item(132, 12).Value = 7382.1432  # L132: 11334.8568 -> 7382.1432
$ws.Cells.Item(132, 13).Value = 286.6774  # M132: -945.3335000000002 -> 286.6774
$ws.Cells.Item(132, 14).Value = -12442.1432  # N132: -16394.8568 -> -12442.1432
$ws.Cells.Item(136, 8).Value = 3789.7837  # H136: 5229.577 -> 3789.7837
$ws.Cells.Item(136, 9).Value = 1193.5454  # I136: 1681.4286 -> 1193.5454
$ws.Cells.Item(136, 10).Value = 7597.6  # J136: 9369.083000000001 -> 7597.6
$ws.Cells.Item(136, 11).Value = 3580.6362  # K136: 5044.2858 -> 3580.6362
$ws.Cells.Item(136, 12).Value = 22792.8  # L136: 28107.249 -> 22792.8
$ws.Cells.Item(136, 13).Value = -1030.6362  # M136: -2494.2858 -> -1030.6362
$ws.Cells.Item(136, 14).Value = -27892.8  # N136: -33207.249 -> -27892.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 14680  # H87: 13800 -> 14680
$ws.Cells.Item(87, 9).Value = 4175  # I87: 5950 -> 4175
$ws.Cells.Item(87, 10).Value = 21683.334  # J87: 21650 -> 21683.334
$ws.Cells.Item(87, 11).Value = 12525  # K87: 17850 -> 12525
$ws.Cells.Item(87, 12).Value = 65050.00199999999  # L87: 64950 -> 65050.00199999999
$ws.Cells.Item(87, 13).Value = -11277  # M87: -16602 -> -11277
$ws.Cells.Item(87, 14).Value = -67546.00199999999  # N87: -67446 -> -67546.00199999999
$ws.Cells.Item(90, 8).Value = 14680  # H90: 13800 -> 14680
$ws.Cells.Item(90, 9).Value = 4175  # I90: 5950 -> 4175
$ws.Cells.Item(90, 10).Value = 21683.334  # J90: 21650 -> 21683.334
$ws.Cells.Item(90, 11).Value = 37575  # K90: 53550 -> 37575
$ws.Cells.Item(90, 12).Value = 195150.006  # L90: 194850 -> 195150.006
$ws.Cells.Item(90, 13).Value = -31335  # M90: -47310 -> -31335
$ws.Cells.Item(90, 14).Value = -207630.006  # N90: -207330 -> -207630.006
$ws.Cells.Item(107, 8).Value = 458429.72  # H107: 519335.88 -> 458429.72
$ws.Cells.Item(107, 10).Value = 1297190.6  # J107: 1944969.2 -> 1297190.6
$ws.Cells.Item(107, 12).Value = 3891571.8  # L107: 5834907.6 -> 3891571.8
$ws.Cells.Item(107, 14).Value = -3895411.8  # N107: -5838747.6 -> -3895411.8
$ws.Cells.Item(117, 8).Value = 2141.8635  # H117: 2142.318 -> 2141.8635
$ws.Cells.Item(117, 10).Value = 2767.875  # J117: 2768.5 -> 2767.875
$ws.Cells.Item(117, 12).Value = 8303.625  # L117: 8305.5 -> 8303.625
$ws.Cells.Item(117, 14).Value = -15187.625  # N117: -15189.5 -> -15187.625
$ws.Cells.Item(129, 8).Value = 41476.52  # H129: 43183.043 -> 41476.52
$ws.Cells.Item(129, 10).Value = 64018.312  # J129: 68251.53 -> 64018.312
$ws.Cells.Item(129, 12).Value = 192054.936  # L129: 204754.59 -> 192054.936
$ws.Cells.Item(129, 14).Value = -202054.936  # N129: -214754.59 -> -202054.936
$ws.Cells.Item(131, 8).Value = 787.3125  # H131: 788.8 -> 787.3125
$ws.Cells.Item(131, 10).Value = 811.8  # J131: 812.34045 -> 811.8
$ws.Cells.Item(131, 12).Value = 2435.4  # L131: 2437.02135 -> 2435.4
$ws.Cells.Item(131, 14).Value = -12515.4  # N131: -12517.02135 -> -12515.4
$ws.Cells.Item(133, 8).Value = 800  # H133: 10000 -> 800
$ws.Cells.Item(133, 9).Value = 800  # I133: 0 -> 800
$ws.Cells.Item(133, 10).Value = 0  # J133: 10000 -> 0
$ws.Cells.Item(133, 11).Value = 2400  # K133: 0 -> 2400
$ws.Cells.Item(133, 12).Value = 0  # L133: 30000 -> 0
$ws.Cells.Item(133, 13).Value = 2660  # M133: None -> 2660
$ws.Cells.Item(133, 14).ClearContents()  # N133: was -40120, now blank

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 4818  # H33: 9990 -> 4818
$ws.Cells.Item(33, 10).Value = 4818  # J33: 9990 -> 4818
$ws.Cells.Item(33, 12).Value = 4818  # L33: 9990 -> 4818
$ws.Cells.Item(33, 14).Value = -5322  # N33: -10494 -> -5322
$ws.Cells.Item(48, 8).Value = 0  # H48: 5000 -> 0
$ws.Cells.Item(48, 9).Value = 0  # I48: 5000 -> 0
$ws.Cells.Item(48, 11).Value = 0  # K48: 5000 -> 0
$ws.Cells.Item(48, 13).ClearContents()  # M48: was -4515, now blank
$ws.Cells.Item(70, 8).Value = 5058.5  # H70: 4541.6665 -> 5058.5
$ws.Cells.Item(70, 9).Value = 3411.3333  # I70: 3450 -> 3411.3333
$ws.Cells.Item(70, 11).Value = 3411.3333  # K70: 3450 -> 3411.3333
$ws.Cells.Item(70, 13).Value = -3141.3333  # M70: -3180 -> -3141.3333
$ws.Cells.Item(73, 8).Value = 5058.5  # H73: 4541.6665 -> 5058.5
$ws.Cells.Item(73, 9).Value = 3411.3333  # I73: 3450 -> 3411.3333
$ws.Cells.Item(73, 11).Value = 3411.3333  # K73: 3450 -> 3411.3333
$ws.Cells.Item(73, 13).Value = -2475.3333  # M73: -2514 -> -2475.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 80965.21000000001  # H100: 1869.8 -> 80965.21000000001
$ws.Cells.Item(100, 9).Value = 112741.3  # I100: 1924.75 -> 112741.3
$ws.Cells.Item(100, 10).Value = 1525  # J100: 1650 -> 1525
$ws.Cells.Item(100, 11).Value = 112741.3  # K100: 1924.75 -> 112741.3
$ws.Cells.Item(100, 12).Value = 1525  # L100: 1650 -> 1525
$ws.Cells.Item(100, 13).Value = -112200.3  # M100: -1383.75 -> -112200.3
$ws.Cells.Item(100, 14).Value = -2607  # N100: -2732 -> -2607
$ws.Cells.Item(132, 8).Value = 1665.8  # H132: 1872.9032 -> 1665.8
$ws.Cells.Item(132, 9).Value = 1268.3077  # I132: 1383.9048 -> 1268.3077
$ws.Cells.Item(132, 10).Value = 4249.5  # J132: 2899.8 -> 4249.5
$ws.Cells.Item(132, 11).Value = 3804.9231  # K132: 4151.7144 -> 3804.9231
$ws.Cells.Item(132, 12).Value = 12748.5  # L132: 8699.400000000001 -> 12748.5
$ws.Cells.Item(132, 13).Value = -1274.9231  # M132: -1621.7144 -> -1274.9231
$ws.Cells.Item(132, 14).Value = -17808.5  # N132: -13759.4 -> -17808.5
$ws.Cells.Item(136, 8).Value = 3919.2856  # H136: 3481.875 -> 3919.2856
$ws.Cells.Item(136, 9).Value = 1149.6  # I136: 1174 -> 1149.6
$ws.Cells.Item(136, 10).Value = 27000  # J136: 7328.3335 -> 27000
$ws.Cells.Item(136, 11).Value = 3448.8  # K136: 3522 -> 3448.8
$ws.Cells.Item(136, 12).Value = 81000  # L136: 21985.0005 -> 81000
$ws.Cells.Item(136, 13).Value = -898.7999999999997  # M136: -972 -> -898.7999999999997
$ws.Cells.Item(136, 14).Value = -86100  # N136: -27085.0005 -> -86100
$ws.Cells.Item(137, 8).Value = 34143  # H137: 0 -> 34143
$ws.Cells.Item(137, 10).Value = 34143  # J137: 0 -> 34143
$ws.Cells.Item(137, 12).Value = 34143  # L137: 0 -> 34143
$ws.Cells.Item(137, 14).Value = -44343  # N137: None -> -44343
$ws.Cells.Item(138, 8).Value = 39433.332  # H138: 39880 -> 39433.332
$ws.Cells.Item(138, 10).Value = 39433.332  # J138: 39880 -> 39433.332
$ws.Cells.Item(138, 12).Value = 39433.332  # L138: 39880 -> 39433.332
$ws.Cells.Item(138, 14).Value = -49713.332  # N138: -50160 -> -49713.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(48, 8).Value = 0  # H48: 1065 -> 0
$ws.Cells.Item(48, 10).Value = 0  # J48: 1065 -> 0
$ws.Cells.Item(48, 12).Value = 0  # L48: 1065 -> 0
$ws.Cells.Item(48, 14).ClearContents()  # N48: was -2203, now blank
$ws.Cells.Item(132, 8).Value = 1608.6  # H132: 1637.4688 -> 1608.6
$ws.Cells.Item(132, 9).Value = 999.9048  # I132: 1015.8421 -> 999.9048
$ws.Cells.Item(132, 10).Value = 2521.6428  # J132: 2546 -> 2521.6428
$ws.Cells.Item(132, 11).Value = 2999.7144  # K132: 3047.5263 -> 2999.7144
$ws.Cells.Item(132, 12).Value = 7564.928400000001  # L132: 7638 -> 7564.928400000001
$ws.Cells.Item(132, 13).Value = -469.7143999999998  # M132: -517.5263 -> -469.7143999999998
$ws.Cells.Item(132, 14).Value = -12624.9284  # N132: -12698 -> -12624.9284
$ws.Cells.Item(133, 8).Value = 73183  # H133: 76553.75 -> 73183
$ws.Cells.Item(133, 10).Value = 73183  # J133: 76553.75 -> 73183
$ws.Cells.Item(133, 12).Value = 73183  # L133: 76553.75 -> 73183
$ws.Cells.Item(133, 14).Value = -83303  # N133: -86673.75 -> -83303
